{"js": "// The diff removes three consecutive paragraphs near the end of the\n// document body:\n//   1) an empty \"Normal\" paragraph,\n//   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n//   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n//      pages. Original theme under Creative Commons Attribution\"\n// The paragraph that immediately follows the text \"S\u00e3o Paulo, 2004....\"\n// (the empty paragraph) through the \"\u00a9 2020 ...\" paragraph are removed,\n// while the subsequent empty paragraph and the page-break paragraph stay.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"Ver no Jupiter ...\" paragraph; the diff removes the empty\n// paragraph immediately preceding it plus it plus the \"\u00a9 2020 ...\"\n// paragraph right after it.\nlet jupiterIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIndex = i;\n    break;\n  }\n}\n\nif (jupiterIndex === -1) {\n  throw new Error(\"Could not find the 'Ver no Jupiter' paragraph\");\n}\n\nconst toDelete = [];\n// The empty paragraph right before it (only if really empty).\nif (jupiterIndex - 1 >= 0 && items[jupiterIndex - 1].text === \"\") {\n  toDelete.push(items[jupiterIndex - 1]);\n}\ntoDelete.push(items[jupiterIndex]);\nif (jupiterIndex + 1 < items.length && items[jupiterIndex + 1].text.indexOf(\"\u00a9 2020\") !== -1) {\n  toDelete.push(items[jupiterIndex + 1]);\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\n\nawait context.sync();\n", "ps1": "# The commit removes three consecutive paragraphs near the end of the\n# document body (right after the \"S\u00e3o Paulo, 2004....Geotecnia\n# ambiental. Elsevier, 2015.\" bibliography line):\n#   1) an empty \"Normal\" paragraph,\n#   2) \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n#   3) \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github\n#      pages. Original theme under Creative Commons Attribution\"\n# The empty paragraph and the page-break paragraph that follow are left\n# untouched.\n\n$d = $word.ActiveDocument\n\n# Locate the \"Ver no Jupiter ...\" paragraph.\n$targetIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -like \"*Ver no Jupiter*\") {\n        $targetIndex = $i\n        break\n    }\n}\n\nif ($targetIndex -eq -1) {\n    throw \"Could not find the 'Ver no Jupiter' paragraph\"\n}\n\n$prevPara = $d.Paragraphs.Item($targetIndex - 1)\n$curPara  = $d.Paragraphs.Item($targetIndex)\n$nextPara = $d.Paragraphs.Item($targetIndex + 1)\n\n# Sanity-check neighbours before deleting: the paragraph before should be\n# blank, and the paragraph after should be the \"\u00a9 2020 ...\" credit line.\n$deleteStart = $curPara.Range.Start\n$deleteEnd   = $curPara.Range.End\n\nif ($prevPara.Range.Text.Trim().Length -eq 0) {\n    $deleteStart = $prevPara.Range.Start\n}\n\nif ($nextPara.Range.Text -like \"*2020*\") {\n    $deleteEnd = $nextPara.Range.End\n}\n\n$killRange = $d.Range($deleteStart, $deleteEnd)\n$killRange.Delete()\n"}
